$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.47"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'23.02"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'3.608"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'5.400"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'0.05912"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'3.456"
$ws.Range("G7").Value = "'12"
$ws.Range("D8").Value = "'6.534"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.8111"
$ws.Range("G9").Value = "'12"
$ws.Range("D10").Value = "'0.9171"
$ws.Range("G10").Value = "'12"
$ws.Range("B11").Value = "'One"
$ws.Range("C11").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D11").Value = "'0.0005944"
$ws.Range("E11").Value = "'10OneONE"
$ws.Range("G11").Value = "'12"
$ws.Range("B12").Value = "'WazirX"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1420"
$ws.Range("E12").Value = "'11WazirXWRX"
$ws.Range("G12").Value = "'12"
$ws.Range("B13").Value = "'MandalaExchangeToken"
$ws.Range("C13").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.07432"
$ws.Range("E13").Value = "'12MandalaExchangeTokenMDX"
$ws.Range("G13").Value = "'12"
$ws.Range("B14").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("C14").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D14").Value = "'0.03278"
$ws.Range("E14").Value = "'13LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G14").Value = "'12"
$ws.Range("B15").Value = "'ProBitToken"
$ws.Range("C15").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D15").Value = "'0.1310"
$ws.Range("E15").Value = "'14ProBitTokenPROB"
$ws.Range("G15").Value = "'12"
$ws.Range("B16").Value = "'BitrueCoin"
$ws.Range("C16").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D16").Value = "'0.03064"
$ws.Range("E16").Value = "'15BitrueCoinBTR"
$ws.Range("G16").Value = "'12"
$ws.Range("B17").Value = "'BitMartToken"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D17").Value = "'0.09345"
$ws.Range("E17").Value = "'16BitMartTokenBMX"
$ws.Range("G17").Value = "'12"
$ws.Range("B18").Value = "'MCDex"
$ws.Range("C18").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D18").Value = "'3.866"
$ws.Range("E18").Value = "'17MCDexMCB"
$ws.Range("G18").Value = "'12"
$ws.Range("B19").Value = "'BitForexToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D19").Value = "'0.001559"
$ws.Range("E19").Value = "'18BitForexTokenBF"
$ws.Range("G19").Value = "'12"
$ws.Range("B20").Value = "'CoinExToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D20").Value = "'0.04668"
$ws.Range("E20").Value = "'19CoinExTokenCET"
$ws.Range("G20").Value = "'12"
$ws.Range("B21").Value = "'TigerCash"
$ws.Range("C21").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D21").Value = "'0.005976"
$ws.Range("E21").Value = "'20TigerCashTCH"
$ws.Range("G21").Value = "'12"
$ws.Range("B22").Value = "'HotbitToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D22").Value = "'0.004946"
$ws.Range("E22").Value = "'21HotbitTokenHTB"
$ws.Range("G22").Value = "'12"
$ws.Range("B23").Value = "'BitKan"
$ws.Range("C23").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "'0.0009874"
$ws.Range("E23").Value = "'22BitKanKAN"
$ws.Range("G23").Value = "'12"
$ws.Range("B24").Value = "'NitroEx"
$ws.Range("C24").Value = "'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D24").Value = "'0.00009606"
$ws.Range("E24").Value = "'23NitroExNTX"
$ws.Range("G24").Value = "'12"
$ws.Range("B25").Value = "'BTSEToken"
$ws.Range("C25").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.151"
$ws.Range("E25").Value = "'24BTSETokenBTSE"
$ws.Range("G25").Value = "'12"
$ws.Range("B26").Value = "'BitpandaEcosystemToken"
$ws.Range("C26").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3241"
$ws.Range("E26").Value = "'25BitpandaEcosystemTokenBEST"
$ws.Range("G26").Value = "'12"
$ws.Range("D27").Value = "'0.0002902"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03970"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.006175"
$ws.Range("E41").Value = "'40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.1077"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.003002"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.009003"
$ws.Range("E44").Value = "'43LocalTradersLCT"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005194"
$ws.Range("G45").Value = "'12"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.8696"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.002282"
$ws.Range("G48").Value = "'12"
$ws.Range("G49").Value = "'12"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
